$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Marking" row (row 11): right-answer marks from 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): total correct marks from 42 -> 70 (14 * 5)
$ws.Range("B12").Value = 70

# Update corr/total text from "35/84" -> "70/140"
$ws.Range("E12").Value = "70/140"
